$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(176, 1).Value = 237675826419
$ws.Cells.Item(176, 2).Value = 'Sara Gaetane Njouma epse Ntoma'
$ws.Cells.Item(176, 3).Value = 'Rte_0'
$ws.Cells.Item(176, 4).Value = 'Esg Building'
$ws.Cells.Item(176, 5).Value = 40765
$ws.Cells.Item(176, 6).Value = 27986
$ws.Cells.Item(176, 7).Value = -12779
$ws.Cells.Item(176, 8).Value = 0.68652029927634
$ws.Cells.Item(176, 9).Value = 'Ndogbong'

$ws.Cells.Item(177, 1).Value = 237676286294
$ws.Cells.Item(177, 2).Value = 'TECLAIRE KAMENI TCHOKOTHE KEMAYOU EPSE TCHANI'
$ws.Cells.Item(177, 3).Value = 'Rte_4'
$ws.Cells.Item(177, 4).Value = 'Esg Building'
$ws.Cells.Item(177, 5).Value = 43870
$ws.Cells.Item(177, 6).Value = 236900
$ws.Cells.Item(177, 7).Value = 193030
$ws.Cells.Item(177, 8).Value = 5.400045589240939
$ws.Cells.Item(177, 9).Value = 'Ndogbong'

$ws.Cells.Item(178, 1).Value = 237677313421
$ws.Cells.Item(178, 2).Value = 'N A PATEM CLOVIS LAHVET'
$ws.Cells.Item(178, 3).Value = 'Rte_4'
$ws.Cells.Item(178, 4).Value = 'Esg Building'
$ws.Cells.Item(178, 5).Value = 233104.6
$ws.Cells.Item(178, 6).Value = 2129010
$ws.Cells.Item(178, 7).Value = 1895905.4
$ws.Cells.Item(178, 8).Value = 9.133281797098812
$ws.Cells.Item(178, 9).Value = 'Ndogbong'

$ws.Cells.Item(179, 1).Value = 237677880357
$ws.Cells.Item(179, 2).Value = 'vanissa sandjouon'
$ws.Cells.Item(179, 3).Value = 'Rte_4'
$ws.Cells.Item(179, 4).Value = 'Esg Building'
$ws.Cells.Item(179, 5).Value = 149646.6666666667
$ws.Cells.Item(179, 6).Value = 19657
$ws.Cells.Item(179, 7).Value = -129989.6666666667
$ws.Cells.Item(179, 8).Value = 0.1313560832182474
$ws.Cells.Item(179, 9).Value = 'Ndogbong'

$ws.Cells.Item(180, 1).Value = 237678201584
$ws.Cells.Item(180, 2).Value = 'Fon Hortencia Engochan ONYEKA CONNECTION'
$ws.Cells.Item(180, 3).Value = 'Rte_0'
$ws.Cells.Item(180, 4).Value = 'Esg Building'
$ws.Cells.Item(180, 5).Value = 500000
$ws.Cells.Item(180, 6).Value = 2866981
$ws.Cells.Item(180, 7).Value = 2366981
$ws.Cells.Item(180, 8).Value = 5.733962
$ws.Cells.Item(180, 9).Value = 'Ndogbong'

$ws.Cells.Item(181, 1).Value = 237678239927
$ws.Cells.Item(181, 2).Value = 'LA NEGRESSE SARL MADAGAL MOHAMED'
$ws.Cells.Item(181, 3).Value = 'Rte_0'
$ws.Cells.Item(181, 4).Value = 'Esg Building'
$ws.Cells.Item(181, 5).Value = 12360.93333333334
$ws.Cells.Item(181, 6).Value = 2337
$ws.Cells.Item(181, 7).Value = -10023.93333333334
$ws.Cells.Item(181, 8).Value = 0.189063393271274
$ws.Cells.Item(181, 9).Value = 'Ndogbong'

$ws.Cells.Item(182, 1).Value = 237678623874
$ws.Cells.Item(182, 2).Value = 'BOUDIEU RHOTA KAMILAH CONNECTION'
$ws.Cells.Item(182, 3).Value = 'Rte_4'
$ws.Cells.Item(182, 4).Value = 'Esg Building'
$ws.Cells.Item(182, 5).Value = 94100
$ws.Cells.Item(182, 6).Value = 180462
$ws.Cells.Item(182, 7).Value = 86362
$ws.Cells.Item(182, 8).Value = 1.917768331562168
$ws.Cells.Item(182, 9).Value = 'Ndogbong'

Write-Output "Added rows 176-182"